$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last student row (row 11) into the new
# row (row 12) so the new row gets the same styles (borders, fonts, text
# number format, etc.) without disturbing the shared style table.
$ws.Range("A11:V11").Copy()
$ws.Range("A12:V12").PasteSpecial(-4122)

# Fill in the new student's data for row 12.
$ws.Range("A12").Value = "200852"
$ws.Range("B12").Value = "Susanne Streng"
$ws.Range("C12").Value = "RS"
$ws.Range("E12").Value = "3+"
$ws.Range("F12").Value = "4"
$ws.Range("G12").Value = "/"
$ws.Range("H12").Value = "3"
$ws.Range("I12").Value = "2"
$ws.Range("J12").Value = "3"
$ws.Range("K12").Value = "3"
$ws.Range("L12").Value = "3"
$ws.Range("M12").Value = "5"
$ws.Range("N12").Value = "5"
$ws.Range("O12").Value = "3"
$ws.Range("P12").Value = "3"
$ws.Range("Q12").Value = "2"
$ws.Range("R12").Value = "2"
$ws.Range("S12").Value = "3"
$ws.Range("T12").Value = "2"
$ws.Range("U12").Value = "2"
$ws.Range("V12").Value = "2"

# Mirror the workbook's recorded selection/active cell (moved to B13 after
# entering the new row of data).
$ws.Range("B13").Select()
